$wb = $excel.ActiveWorkbook

# --- Create "Week 3" sheet as a copy of "Week 2", placed between "Week 2" and "Assignment" ---
$week2 = $wb.Worksheets.Item("Week 2")
$week2.Copy($null, $week2)
$week3 = $wb.Worksheets.Item("Week 2 (2)")
$week3.Name = "Week 3"

# The new sheet was an exact copy of "Week 2" (rows 2-25). Shift everything down
# one row so the new sheet's data occupies rows 3-26, as in the target layout.
$week3.Rows.Item(1).Insert()

# Update the title cell to reference the Week 3 label.
$week3.Range("B3").Value = "Web Devlopment Internship 2024 (Week 3)"

# Update the date header row: only the first day (column D) has a date so far;
# the other three day columns are still blank for this new week.
$week3.Range("D4").Value = 45462
$week3.Range("E4:G4").ClearContents()

# Clear all attendance marks in columns E, F, G (days 2-4 haven't happened yet).
$week3.Range("E5:G26").ClearContents()

# Only a subset of students are marked present ("p") for day 1 (column D).
$week3.Range("D8").Value = "p"
$week3.Range("D11").Value = "p"
$week3.Range("D17").Value = "p"
$week3.Range("D19").Value = "p"
$week3.Range("D20").Value = "p"
$week3.Range("D13").ClearContents()
$week3.Range("D24").ClearContents()
$week3.Range("D26").ClearContents()

# --- "Week 2" is no longer the active tab; update its lingering selection. ---
$week2.Activate()
$week2.Range("B2:H25").Select()

# Make "Week 3" the active/visible tab with its own selection.
$week3.Activate()
$excel.ActiveWindow.Zoom = 130
$week3.Range("G9").Select()
